$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I38").Value = 931.44446
$ws.Range("M38").Value = -2422.33338
$ws.Range("K38").Value = 2794.33338
$ws.Range("H38").Value = 931.44446
$ws.Range("K100").Value = 2133.3333
$ws.Range("I100").Value = 2133.3333
$ws.Range("H100").Value = 1902.5
$ws.Range("M100").Value = -1592.3333
$ws.Range("N134").Value = -130140
$ws.Range("H134").Value = 120000
$ws.Range("L134").Value = 120000
$ws.Range("J134").Value = 120000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9042.703
$ws.Range("K32").Value = 6130.6
$ws.Range("L32").Value = 12468.706
$ws.Range("M32").Value = -5843.6
$ws.Range("N32").Value = -13042.706
$ws.Range("I32").Value = 6130.6
$ws.Range("J32").Value = 12468.706
$ws.Range("N45").Value = -2744
$ws.Range("J45").Value = 1990
$ws.Range("H45").Value = 1990
$ws.Range("L45").Value = 1990
$ws.Range("J63").Value = 3133
$ws.Range("H63").Value = 2722.625
$ws.Range("N63").Value = -4505
$ws.Range("L63").Value = 3133
$ws.Range("J66").Value = 3133
$ws.Range("H66").Value = 2722.625
$ws.Range("L66").Value = 15665
$ws.Range("N66").Value = -22529
$ws.Range("I74").Value = 17859928
$ws.Range("M74").Value = -17859054
$ws.Range("H74").Value = 17859928
$ws.Range("K74").Value = 17859928
$ws.Range("M77").Value = -89295272
$ws.Range("I77").Value = 17859928
$ws.Range("K77").Value = 89299640
$ws.Range("H77").Value = 17859928
$ws.Range("M97").Value = 71.88889999999998
$ws.Range("H97").Value = 451.54544
$ws.Range("N97").Value = -1567
$ws.Range("K97").Value = 424.1111
$ws.Range("J97").Value = 575
$ws.Range("I97").Value = 424.1111
$ws.Range("L97").Value = 575
$ws.Range("K110").Value = 127860.125
$ws.Range("H110").Value = 103238.9
$ws.Range("I110").Value = 127860.125
$ws.Range("M110").Value = -125815.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 395998
$ws.Range("J70").Value = 395998
$ws.Range("L70").Value = 395998
$ws.Range("N70").Value = -396584
$ws.Range("L73").Value = 395998
$ws.Range("H73").Value = 395998
$ws.Range("J73").Value = 395998
$ws.Range("N73").Value = -398026
$ws.Range("L94").Value = 1000
$ws.Range("J94").Value = 1000
$ws.Range("M94").Value = -269
$ws.Range("I94").Value = 720
$ws.Range("H94").Value = 776
$ws.Range("N94").Value = -1902
$ws.Range("K94").Value = 720
$ws.Range("L133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("H133").Value = 0
$ws.Range("M134").Value = -88243143
$ws.Range("N134").Value = -18468
$ws.Range("H134").Value = 22730962
$ws.Range("L134").Value = 13398
$ws.Range("J134").Value = 4466
$ws.Range("K134").Value = 88245678
$ws.Range("I134").Value = 29415226
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("K7").Value = 7724
$ws.Range("N7").Value = -312.2
$ws.Range("I7").Value = 7724
$ws.Range("H7").Value = 5602.3887
$ws.Range("L7").Value = 86.2
$ws.Range("M7").Value = -7611
$ws.Range("J7").Value = 86.2
$ws.Range("M86").Value = -3055.2856
$ws.Range("K86").Value = 4178.2856
$ws.Range("H86").Value = 4031
$ws.Range("I86").Value = 4178.2856
$ws.Range("M89").Value = -15275.428
$ws.Range("I89").Value = 4178.2856
$ws.Range("K89").Value = 20891.428
$ws.Range("H89").Value = 4031
$ws.Range("H125").Value = 59498
$ws.Range("L125").Value = 59498
$ws.Range("N125").Value = -64418
$ws.Range("J125").Value = 59498
$ws.Range("J132").Value = 600
$ws.Range("I132").Value = 71429464
$ws.Range("K132").Value = 214288392
$ws.Range("M132").Value = -214285862
$ws.Range("L132").Value = 1800
$ws.Range("N132").Value = -6860
$ws.Range("N138").Value = -125280
$ws.Range("H138").Value = 115000
$ws.Range("L138").Value = 115000
$ws.Range("J138").Value = 115000

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1448.5834
$ws.Range("L107").Value = 5286.4998
$ws.Range("N107").Value = -9126.4998
$ws.Range("J107").Value = 1762.1666
$ws.Range("H125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("M136").Value = 2403
$ws.Range("K136").Value = 2697
$ws.Range("I136").Value = 899
$ws.Range("H136").Value = 899
$ws.Range("H138").Value = 4143.6
$ws.Range("K138").Value = 12430.8
$ws.Range("I138").Value = 4143.6
$ws.Range("M138").Value = -7290.800000000001
$ws.Range("N125").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M97").Value = -1167.2222
$ws.Range("H97").Value = 1924.3636
$ws.Range("N97").Value = -4091.5
$ws.Range("K97").Value = 1663.2222
$ws.Range("J97").Value = 3099.5
$ws.Range("I97").Value = 1663.2222
$ws.Range("L97").Value = 3099.5
$ws.Range("K107").Value = 891.17645
$ws.Range("I107").Value = 891.17645
$ws.Range("H107").Value = 1238.32
$ws.Range("M107").Value = 1028.82355
$ws.Range("H121").Value = 14975
$ws.Range("L121").Value = 14975
$ws.Range("N121").Value = -18469
$ws.Range("J121").Value = 14975

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K16").Value = 965.2727
$ws.Range("M16").Value = -795.2727
$ws.Range("I16").Value = 965.2727
$ws.Range("H16").Value = 2800.3333
$ws.Range("H61").Value = 2451.1333
$ws.Range("I61").Value = 2004.36
$ws.Range("M61").Value = -1802.36
$ws.Range("K61").Value = 2004.36
$ws.Range("M68").Value = -2315214.8
$ws.Range("K68").Value = 2315963.8
$ws.Range("H68").Value = 1895606.6
$ws.Range("N68").Value = -5498
$ws.Range("I68").Value = 2315963.8
$ws.Range("L68").Value = 4000
$ws.Range("J68").Value = 4000
$ws.Range("N71").Value = -27488
$ws.Range("K71").Value = 11579819
$ws.Range("M71").Value = -11576075
$ws.Range("J71").Value = 4000
$ws.Range("H71").Value = 1895606.6
$ws.Range("L71").Value = 20000
$ws.Range("I71").Value = 2315963.8
$ws.Range("I93").Value = 1034
$ws.Range("K93").Value = 1034
$ws.Range("M93").Value = 214
$ws.Range("H93").Value = 1034
$ws.Range("K113").Value = 2004.36
$ws.Range("H113").Value = 2451.1333
$ws.Range("I113").Value = 2004.36
$ws.Range("M113").Value = 165.6400000000001
$ws.Range("H117").Value = 99392
$ws.Range("L117").Value = 99392
$ws.Range("N117").Value = -108570
$ws.Range("J117").Value = 99392
$ws.Range("L133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("H133").Value = 0
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("L46").Value = 63494.75
$ws.Range("J46").Value = 63494.75
$ws.Range("H46").Value = 63494.75
$ws.Range("N46").Value = -63956.75
$ws.Range("I81").Value = 2849.5
$ws.Range("L81").Value = 4332.6666
$ws.Range("M81").Value = -4638
$ws.Range("N81").Value = -6454.6666
$ws.Range("H81").Value = 2439.6
$ws.Range("J81").Value = 2166.3333
$ws.Range("K81").Value = 5699
$ws.Range("H84").Value = 2439.6
$ws.Range("L84").Value = 21663.333
$ws.Range("J84").Value = 2166.3333
$ws.Range("I84").Value = 2849.5
$ws.Range("M84").Value = -23191
$ws.Range("K84").Value = 28495
$ws.Range("N84").Value = -32271.333
$ws.Range("L87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("H87").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("M96").Value = -1724.8
$ws.Range("I96").Value = 3097.8
$ws.Range("K96").Value = 3097.8
$ws.Range("H96").Value = 3936.0356
$ws.Range("N100").Value = -3060
$ws.Range("K100").Value = 2350.8
$ws.Range("I100").Value = 1175.4
$ws.Range("H100").Value = 1138.12
$ws.Range("M100").Value = -1809.8
$ws.Range("J100").Value = 989
$ws.Range("L100").Value = 1978
$ws.Range("H123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("J132").Value = 2441.8333
$ws.Range("I132").Value = 17247710
$ws.Range("H132").Value = 14291379
$ws.Range("K132").Value = 51743130
$ws.Range("M132").Value = -51740600
$ws.Range("L132").Value = 7325.499899999999
$ws.Range("N132").Value = -12385.4999
$ws.Range("N134").Value = -195554.25
$ws.Range("H134").Value = 63494.75
$ws.Range("L134").Value = 190484.25
$ws.Range("J134").Value = 63494.75
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()
$ws.Range("N123").ClearContents()
